# Updated cryptos list on Thu Jul 25 17:12:33 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for every coin row
# with newly-scraped figures, and fixes three rows whose (Coin, Link, Price,
# Volume) tuples had drifted out of their correct ranking position:
#   - rows 39/40 swap back to Stacks / EnergySwap
#   - rows 45/46/47 rotate back to Hedera / InjectiveProtocol / OKB
#
# All of these sheet cells are stored as literal text (prices such as
# "64.881.08" use '.' as a thousands separator, not a decimal point, and the
# volume column keeps its padding spaces), so every write below forces the
# cell to Text format first and clears that temporary formatting afterwards
# -- this stops Excel from "helpfully" re-interpreting e.g. "1.00" as the
# number 1 and dropping the trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Ref,
        [string]$Text
    )
    $range = $ws.Range($Ref)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.ClearFormats()
}

# row, B (Coin), C (Link), D (Price), E (Volume 1h)
$rows = @(
    @{ Row = 2;  D = "65.004.94"; E = "  -2.24%  " },
    @{ Row = 3;  D = "3.166.80";  E = "  -7.50%  " },
    @{ Row = 4;                  E = "  +0.01%  " },
    @{ Row = 5;  D = "569.86";    E = "  -2.65%  " },
    @{ Row = 6;  D = "171.61";    E = "  -5.12%  " },
    @{ Row = 7;  D = "0.623";     E = "  -0.08%  " },
    @{ Row = 8;                  E = "  +0.03%  " },
    @{ Row = 9;  D = "3.164.70";  E = "  -7.43%  " },
    @{ Row = 10;                 E = "  -5.49%  " },
    @{ Row = 11;                 E = "  -5.67%  " },
    @{ Row = 12; D = "0.396";     E = "  -4.44%  " },
    @{ Row = 13; D = "3.710.74";  E = "  -7.56%  " },
    @{ Row = 14;                 E = "  +1.13%  " },
    @{ Row = 15; D = "27.23";     E = "  -6.70%  " },
    @{ Row = 16; D = "64.964.16"; E = "  -2.18%  " },
    @{ Row = 17;                 E = "  -5.61%  " },
    @{ Row = 18; D = "3.158.36";  E = "  -7.12%  " },
    @{ Row = 19; D = "5.74";      E = "  -3.11%  " },
    @{ Row = 20; D = "12.88";     E = "  -7.05%  " },
    @{ Row = 21; D = "357.15";    E = "  -3.47%  " },
    @{ Row = 22; D = "7.28";      E = "  -4.35%  " },
    @{ Row = 23; D = "1.00";      E = "  +0.21%  " },
    @{ Row = 24; D = "69.18";     E = "  -5.66%  " },
    @{ Row = 25; D = "0.501";     E = "  -6.51%  " },
    @{ Row = 26;                 E = "  -7.17%  " },
    @{ Row = 27; D = "9.72";      E = "  -1.34%  " },
    @{ Row = 28;                 E = "  -2.17%  " },
    @{ Row = 29; D = "0.997";     E = "  -0.32%  " },
    @{ Row = 30;                 E = "  -0.15%  " },
    @{ Row = 31; D = "1.92";      E = "  -4.03%  " },
    @{ Row = 32;                 E = "  -7.45%  " },
    @{ Row = 33; D = "22.04";     E = "  -5.81%  " },
    @{ Row = 34; D = "6.67";      E = "  -5.72%  " },
    @{ Row = 35;                 E = "  -4.68%  " },
    @{ Row = 36;                 E = "  -6.22%  " },
    @{ Row = 37; D = "154.35";    E = "  -5.26%  " },
    @{ Row = 38; D = "0.839";     E = "  -3.34%  " },
    @{ Row = 39; B = "Stacks";     C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx";              D = "1.77";  E = "  -1.93%  " },
    @{ Row = 40; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";          D = "26.11"; E = "  -5.49%  " },
    @{ Row = 41;                 E = "  -4.62%  " },
    @{ Row = 42; D = "2.657.96";  E = "  -2.63%  " },
    @{ Row = 43; D = "4.19";      E = "  -5.60%  " },
    @{ Row = 44;                 E = "  -4.19%  " },
    @{ Row = 45; B = "Hedera";            C = "https://coinranking.com/coin/jad286TjB+hedera-hbar";             D = "0.0661"; E = "  -4.13%  " },
    @{ Row = 46; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj";   D = "24.41";  E = "  -2.75%  " },
    @{ Row = 47; B = "OKB";               C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb";             D = "39.42";  E = "  -1.20%  " },
    @{ Row = 48; D = "326.91";    E = "  -2.88%  " },
    @{ Row = 49;                 E = "  -4.17%  " },
    @{ Row = 50; D = "0.103";     E = "  -1.44%  " },
    @{ Row = 51; D = "0.999";     E = "  -0.01%  " }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    if ($entry.ContainsKey("B")) { Set-TextValue ("B" + $r) $entry.B }
    if ($entry.ContainsKey("C")) { Set-TextValue ("C" + $r) $entry.C }
    if ($entry.ContainsKey("D")) { Set-TextValue ("D" + $r) $entry.D }
    if ($entry.ContainsKey("E")) { Set-TextValue ("E" + $r) $entry.E }
}
